$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B, shifting en_list/en_1/en_2 one column right
$ws.Columns("B").Insert() | Out-Null

# New header for the inserted column
$ws.Range("B1").Value = "en_comments"

# Autofit the new column's width (best-effort approximation of Excel's
# font-metric "best fit" width)
$ws.Columns("B").AutoFit() | Out-Null

# Expand the existing color-scale conditional formatting to cover the new column
$fc = $ws.Range("A1").FormatConditions.Item(1)
$fc.ModifyAppliesToRange($ws.Range("A1:B1")) | Out-Null

# Move the active selection to match the saved workbook state
$ws.Range("B11").Select() | Out-Null
